$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.524.56"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.463.96"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("E9").Value = "  +4.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "2.845.86"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "2.461.93"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "41.526.10"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.31%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("E36").Value = "  -7.77%  "
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  -3.64%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "1.942.43"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "2.704.42"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.91%  "
